$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'basketball gear boys',
    'spandex shorts for volleyball',
    'volleyball knee sleeves men',
    'knee pads for work construction',
    'knee pads for working on floors',
    'bee tights adult',
    'compression test',
    'lacrosse youth shorts',
    'women volleyball knee pads',
    'short football pants',
    'girls compression tights',
    'knee compression sleeve reduce strain & swelling',
    'volleyball knee pads small',
    'yoga pants knee length',
    'pad for squats',
    'sweat pads',
    'calf protector',
    'compression pants long',
    'spandex mens pants',
    'protective knee pads for work',
    'youth compression shorts',
    'mens basketball shorts long',
    'thick volleyball knee pads',
    'weightlifting guide',
    'small work knee pads',
    'basketballs in bulk',
    'baseballs cheap',
    'spandex capri leggings',
    'hex fabric',
    'black knee pads for work',
    'black knee guards',
    'work knee pads for men gel',
    'leg guard baseball',
    'lacrosse shorts youth',
    'girdles for men',
    'basketball cheap',
    'wrestling fight shorts',
    'construction knee',
    'knee compression sleeve protector',
    'knee bursitis sleeve',
    'football leggings',
    'black girls softball pants',
    'girls softball pants black',
    'working pants with knee pads',
    'youth volleyball sleeves',
    'mens knee length shorts',
    'compression capris girls',
    'leaf leggings',
    'fit compression knee',
    'performance basketball',
    'pantalones de basketball',
    'youth baseball sleeves for boys',
    'mens tights for sports',
    'knee protector work',
    'long basketball shorts for men',
    'mens knee sleeves weightlifting',
    'knee pads bulk',
    'snowboarding protective gear',
    'calf compression sleeve youth',
    'knee pain pads',
    'knee pad for construction',
    'mens capri yoga pants',
    'knee sleeves for basketball',
    'knee pad sleeves',
    'knee work',
    'volleyball shorts men',
    'knee pad for working',
    'above knee shorts men',
    'youth volleyball shorts for girls',
    'protector paintball',
    'baseball shorts boys',
    'boys tights and leggings',
    'kneeling pad gel',
    'paintball pants men',
    'men above knee shorts',
    'spandex compression shorts',
    'hockey padded shorts',
    'knee pads for joint pain',
    'rodillera volleyball',
    'basketball calf sleeve',
    'womens lacrosse pants',
    'best work knee pads',
    'leg sleeves for men basketball',
    'football pouch youth',
    'calf compression sleeve boys',
    'boys leggings sports',
    'basketball stretch pants',
    'mens spandex shorts',
    'black kneepads',
    'joint protectors',
    'boys hiking pants',
    'knee length yoga pants',
    'girls volleyball shorts youth',
    'baseball pants mens long',
    'mens sports tights',
    'womens football pads',
    'stretch mark men',
    'youth xl baseball pants',
    'youth basketball compression sleeve',
    'sliding knee sleeve'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
